# Logged Week 15 and simulated Week 16
# Update the "R" (row 3) totals on both the OFF and DEF sheets.

$wb = $excel.ActiveWorkbook

# --- OFF sheet ---
$wsOff = $wb.Worksheets.Item("OFF")
$wsOff.Range("B3").Value = 177
$wsOff.Range("C3").Value = 129
$wsOff.Range("D3").Value = 62
$wsOff.Range("E3").Value = 25
$wsOff.Range("G3").Value = 2

# --- DEF sheet ---
$wsDef = $wb.Worksheets.Item("DEF")
$wsDef.Range("B3").Value = 183
$wsDef.Range("C3").Value = 138
$wsDef.Range("D3").Value = 45
$wsDef.Range("E3").Value = 20
$wsDef.Range("G3").Value = 4
